$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1199.8125
$ws.Range("J17").Value = 1231.1333
$ws.Range("L17").Value = 3693.3999
$ws.Range("N17").Value = -4029.3999
$ws.Range("H40").Value = 2316.389
$ws.Range("I40").Value = 1517.4166
$ws.Range("K40").Value = 1517.4166
$ws.Range("M40").Value = -1342.4166
$ws.Range("H61").Value = 137
$ws.Range("I61").Value = 137
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 411
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -239
$ws.Range("H64").Value = 3317.6177
$ws.Range("I64").Value = 3020
$ws.Range("J64").Value = 3441.625
$ws.Range("K64").Value = 3020
$ws.Range("L64").Value = 3441.625
$ws.Range("M64").Value = -2772
$ws.Range("N64").Value = -3937.625
$ws.Range("H67").Value = 3317.6177
$ws.Range("I67").Value = 3020
$ws.Range("J67").Value = 3441.625
$ws.Range("K67").Value = 3020
$ws.Range("L67").Value = 3441.625
$ws.Range("M67").Value = -2162
$ws.Range("N67").Value = -5157.625
$ws.Range("H137").Value = 2022.9131
$ws.Range("I137").Value = 2945.3333
$ws.Range("J137").Value = 1016.63635
$ws.Range("K137").Value = 8835.999899999999
$ws.Range("L137").Value = 3049.90905
$ws.Range("M137").Value = -6285.999899999999
$ws.Range("N137").Value = -8149.90905
$ws.Range("H138").Value = 1547.7333
$ws.Range("I138").Value = 1301.2188
$ws.Range("J138").Value = 1829.4642
$ws.Range("K138").Value = 3903.6564
$ws.Range("L138").Value = 5488.392599999999
$ws.Range("M138").Value = 1236.3436
$ws.Range("N138").Value = -15768.3926
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7458.3086
$ws.Range("I32").Value = 6606.8813
$ws.Range("K32").Value = 6606.8813
$ws.Range("M32").Value = -6319.8813
$ws.Range("H122").Value = 894
$ws.Range("I122").Value = 801.2632
$ws.Range("J122").Value = 1481.3334
$ws.Range("K122").Value = 2403.7896
$ws.Range("L122").Value = 4444.0002
$ws.Range("M122").Value = 46.21039999999994
$ws.Range("N122").Value = -9344.0002
$ws.Range("H132").Value = 621599.3
$ws.Range("I132").Value = 1324469.6
$ws.Range("J132").Value = 4444.927
$ws.Range("K132").Value = 3973408.8
$ws.Range("L132").Value = 13334.781
$ws.Range("M132").Value = -3970878.8
$ws.Range("N132").Value = -18394.781
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1674.5
$ws.Range("I86").Value = 1549
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 1549
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -426
$ws.Range("N86").Value = -4046
$ws.Range("H89").Value = 1674.5
$ws.Range("I89").Value = 1549
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 7745
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -2129
$ws.Range("N89").Value = -20232
$ws.Range("H94").Value = 740.1515000000001
$ws.Range("I94").Value = 659.5
$ws.Range("J94").Value = 1039.7142
$ws.Range("K94").Value = 659.5
$ws.Range("L94").Value = 1039.7142
$ws.Range("M94").Value = -208.5
$ws.Range("N94").Value = -1941.7142
$ws.Range("H107").Value = 1822
$ws.Range("I107").Value = 1822
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1822
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 98
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2267.3333
$ws.Range("I99").Value = 1501.7142
$ws.Range("J99").Value = 2937.25
$ws.Range("K99").Value = 1501.7142
$ws.Range("L99").Value = 2937.25
$ws.Range("M99").Value = -3.714199999999892
$ws.Range("N99").Value = -5933.25
$ws.Range("H107").Value = 1569.6471
$ws.Range("I107").Value = 382.7143
$ws.Range("J107").Value = 2400.5
$ws.Range("K107").Value = 382.7143
$ws.Range("L107").Value = 2400.5
$ws.Range("M107").Value = 1537.2857
$ws.Range("N107").Value = -6240.5
$ws.Range("H126").Value = 2267.3333
$ws.Range("I126").Value = 1501.7142
$ws.Range("J126").Value = 2937.25
$ws.Range("K126").Value = 4505.142599999999
$ws.Range("L126").Value = 8811.75
$ws.Range("M126").Value = -2035.142599999999
$ws.Range("N126").Value = -13751.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 513.6875
$ws.Range("I107").Value = 466.125
$ws.Range("J107").Value = 561.25
$ws.Range("K107").Value = 1398.375
$ws.Range("L107").Value = 1683.75
$ws.Range("M107").Value = 521.625
$ws.Range("N107").Value = -5523.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.916668
$ws.Range("I2").Value = 55.555557
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 55.555557
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 57.444443
$ws.Range("N2").Value = -271
$ws.Range("H102").Value = 2219.75
$ws.Range("I102").Value = 1381.4
$ws.Range("J102").Value = 4734.8
$ws.Range("K102").Value = 1381.4
$ws.Range("L102").Value = 4734.8
$ws.Range("M102").Value = 240.5999999999999
$ws.Range("N102").Value = -7978.8
$ws.Range("H113").Value = 1425.5
$ws.Range("I113").Value = 1455.8572
$ws.Range("J113").Value = 1213
$ws.Range("K113").Value = 1455.8572
$ws.Range("L113").Value = 1213
$ws.Range("M113").Value = 714.1428000000001
$ws.Range("N113").Value = -5553
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1526.4706
$ws.Range("I68").Value = 1288.4615
$ws.Range("J68").Value = 2300
$ws.Range("K68").Value = 1288.4615
$ws.Range("L68").Value = 2300
$ws.Range("M68").Value = -539.4614999999999
$ws.Range("N68").Value = -3798
$ws.Range("H71").Value = 1526.4706
$ws.Range("I71").Value = 1288.4615
$ws.Range("J71").Value = 2300
$ws.Range("K71").Value = 6442.307499999999
$ws.Range("L71").Value = 11500
$ws.Range("M71").Value = -2698.307499999999
$ws.Range("N71").Value = -18988
$ws.Range("H82").Value = 1090.6428
$ws.Range("I82").Value = 779.6667
$ws.Range("J82").Value = 1650.4
$ws.Range("K82").Value = 779.6667
$ws.Range("L82").Value = 1650.4
$ws.Range("M82").Value = -418.6667
$ws.Range("N82").Value = -2372.4
$ws.Range("H85").Value = 1090.6428
$ws.Range("I85").Value = 779.6667
$ws.Range("J85").Value = 1650.4
$ws.Range("K85").Value = 779.6667
$ws.Range("L85").Value = 1650.4
$ws.Range("M85").Value = 468.3333
$ws.Range("N85").Value = -4146.4
$ws.Range("H128").Value = 500010900
$ws.Range("J128").Value = 500010900
$ws.Range("L128").Value = 500010900
$ws.Range("N128").Value = -500020860
$ws.Range("H132").Value = 37013.234
$ws.Range("I132").Value = 46626.086
$ws.Range("J132").Value = 5428.143
$ws.Range("K132").Value = 139878.258
$ws.Range("L132").Value = 16284.429
$ws.Range("M132").Value = -137348.258
$ws.Range("N132").Value = -21344.429
